# Updated cryptos list - price/volume refresh and two row content swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "'37.081.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = "'2.058.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'248.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.64%  '
$ws.Range("D6").Value = "'0.657"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'55.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +16.11%  '
$ws.Range("D9").Value = "'61.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.20%  '
$ws.Range("E10").Value = '  +1.06%  '
$ws.Range("D11").Value = "'0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.88%  '
$ws.Range("E12").Value = '  +5.63%  '
$ws.Range("E13").Value = '  +5.37%  '
$ws.Range("D14").Value = "'2.357.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").Value = "'0.816"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").Value = "'2.055.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("D18").Value = "'37.043.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("D19").Value = "'0.0₃0938"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.58%  '
$ws.Range("D20").Value = "'72.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").Value = "'14.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.81%  '
$ws.Range("D22").Value = "'5.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.95%  '
$ws.Range("D23").Value = "'237.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.44%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").Value = "'2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").Value = "'170.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("D28").Value = "'20.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.76%  '
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  +2.28%  '
$ws.Range("E32").Value = '  +11.62%  '
$ws.Range("D33").Value = "'0.0625"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.99%  '
$ws.Range("E34").Value = '  +6.95%  '
$ws.Range("E35").Value = '  +0.06%  '

# --- Rows 36 & 37: Kaspa / LidoDAOToken swapped places (with refreshed price/volume) ---
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = "'2.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.86%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = "'0.0856"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.40%  '

$ws.Range("E38").Value = '  -6.37%  '
$ws.Range("E39").Value = '  +0.63%  '
$ws.Range("D40").Value = "'0.107"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +26.69%  '
$ws.Range("D41").Value = "'18.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.39%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -3.23%  '
$ws.Range("D44").Value = "'96.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("D45").Value = "'4.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +46.14%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("E47").Value = '  +7.00%  '
$ws.Range("D48").Value = "'13.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -51.76%  '
$ws.Range("D49").Value = "'1.297.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.16%  '
$ws.Range("D50").Value = "'2.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.61%  '

# --- Row 51: THORChain replaced by FraxShare ---
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = "'6.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.39%  '
